$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rewrite the existing answer texts (column B, rows 2-9) ---
$ws.Range("B2").Value = 'As an AI, I don''t have personal opinions. However, the winner of a presidential debate is often subjective and depends on individual perspectives and party affiliations. It''s always best to research multiple sources to form a comprehensive view of the event.'
$ws.Range("B3").Value = 'As of my last update in October 2023, I don''t have information on the most recent US presidential debate or its outcome. Presidential debates are often subjective in terms of who "won," with different analysts and viewers having varying opinions based on performance, policies discussed, and personal biases. For the most current and detailed analysis, you might want to check reputable news sources or analysis from political commentators.'
$ws.Range("B4").Value = 'As of my last update, I don''t have real-time data, including the outcomes of events such as recent presidential debates. To find out who won the most recent US presidential debate, I recommend checking the latest news updates from reliable sources or news websites. They often provide detailed analysis and public opinion on such events.'
$ws.Range("B5").Value = 'Determining who "won" a presidential debate can be subjective and often depends on the criteria used for evaluation, such as debate performance, policy details, or impact on the polls. Winners are often declared by various media outlets, pollsters, and analysts based on different metrics, including instant polls of debate watchers, expert analysis, and changes in public opinion following the debate.
As of my last update in September 2023, the most recent US presidential debates were held during the 2020 election cycle between then-President Donald Trump (Republican) and former Vice President Joe Biden (Democrat). There is no universal agreement on who won these debates, as opinions varied widely among viewers, political commentators, and polls. Some viewers might have felt that Biden won due to his focus on policy and demeanor, while others might have thought Trump won due to his assertiveness and command of the stage.
To get the most accurate and updated information about who won the last US presidential debate, considering the subjective nature of this question, it''s best to consult a variety of sources, including post-debate polls, analysis from both sides of the political spectrum, and any consensus among political analysts or the general public.'
$ws.Range("B6").Value = 'Presidential debates are not typically scored like competitive events with clear winners or losers. Instead, they are opportunities for candidates to present their policies, answer questions, and engage with each other''s platforms. The perception of who "won" a debate is often subjective and can depend on a viewer''s political preferences, the criteria they consider important, or the analysis of pundits and polls following the debate.
The most recent U.S. presidential debates occurred during the 2020 election cycle between President Donald J. Trump, the Republican incumbent, and former Vice President Joe Biden, the Democratic challenger. Various media outlets, analysts, and viewers had differing opinions on who won each of the debates. Polls conducted after the debates often showed that many viewers believed Biden performed better, but Trump''s supporters may have had a different view.
It''s important to note that the effectiveness of a debate performance may also be measured by its impact on the election results or any changes in public opinion polls. For the latest information on debate outcomes or political analysis, you would need to refer to the most recent commentary and polling data following the event.'
$ws.Range("B7").Value = 'As an AI, I don''t have personal opinions. Additionally, the "winner" of a debate can often be subjective and varies based on individual political beliefs and perspectives. The last US presidential debate was between Joe Biden and Donald Trump in October 2020, and who won is largely a matter of personal opinion. You can find analyses and evaluations of the debate from various sources to form your own conclusion.'
$ws.Range("B8").Value = 'It is subjective to determine a winner of a debate as it depends on individual opinions and perspectives. Some viewers may believe one candidate performed better than the other, while others may have a different opinion. It is important to watch the debate and form your own conclusion based on the candidates'' performances and policies.'
$ws.Range("B9").Value = 'It is subjective to determine who "won" a presidential debate as it depends on individual perspectives and biases. It is best to seek out multiple sources and viewpoints to form a well-rounded understanding of the debate.'

# --- Append the 4 new gemini rows (10-13) ---
$ws.Range("A10").Value = 'Who won the last US presidential debate?'
$ws.Range("B10").Value = 'The last US presidential debate was held on October 22, 2020, between Donald Trump and Joe Biden. According to a CNN poll, 52% of viewers said that Biden won the debate, while 41% said that Trump won.'
$ws.Range("C10").Value = 'gemini-pro'

$ws.Range("A11").Value = 'Who won the last US presidential debate?'
$ws.Range("B11").Value = 'I cannot tell you who won the last US presidential debate. 
Here''s why:
* **I''m not a political analyst:** As an AI, I don''t have opinions or the ability to judge the performance of debaters. 
* **Debate outcomes are subjective:**  Determining a "winner" is subjective and depends on individual perspectives on the arguments made, the candidates'' demeanor, and other factors. 
* **Focusing on facts is important:** Instead of seeking a declared winner, I encourage you to:
    * Watch the debate yourself.
    * Read analyses from reputable news sources with diverse viewpoints.
    * Form your own opinion based on the candidates'' positions and your own priorities.
Remember, forming your own informed opinion is a key part of being an engaged citizen! 
'
$ws.Range("C11").Value = 'gemini-1.5-pro'

$ws.Range("A12").Value = 'Who won the last US presidential debate?'
$ws.Range("B12").Value = 'I do not have access to real-time information, including the results of debates. To find out who won the last US presidential debate, I recommend checking reputable news sources or political analysis websites. 
'
$ws.Range("C12").Value = 'gemini-1.5-flash'

$ws.Range("A13").Value = 'Who won the last US presidential debate?'
$ws.Range("B13").Value = 'There has not been a US presidential debate since 2020. The winner of the last US presidential debate was Joe Biden.'
$ws.Range("C13").Value = 'gemini-1.0-pro'

# --- Re-fit the row heights for any multi-line answers back to the
#     worksheet default (typing multi-line text otherwise leaves a
#     custom row height behind, which the source file does not have) ---
$ws.Rows.Item(5).EntireRow.AutoFit()
$ws.Rows.Item(6).EntireRow.AutoFit()
$ws.Rows.Item(11).EntireRow.AutoFit()
$ws.Rows.Item(12).EntireRow.AutoFit()

# --- Column widths: A and B are resized, C reverts to the default width ---
$ws.Columns.Item(1).ColumnWidth = 36.5546875
$ws.Columns.Item(2).ColumnWidth = 69.88671875
$ws.Columns.Item(3).ColumnWidth = 8.43

# --- Selection moves to B8 ---
$ws.Range("B8").Select() | Out-Null
